$wb = $excel.ActiveWorkbook

# --- Rename sheets: "Datos" -> "Data", "Ficha técnica" -> "Metadata" ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "Data"
$ws2.Name = "Metadata"

# --- Sheet "Data": sort the year rows (A2:C15) descending by year (col A) ---
$dataRange = $ws1.Range("A2:C15")
$dataRange.Sort($ws1.Range("A2:A15"), 2)

# --- Sheet "Metadata": reorder / relabel the indicator fields, lower-case the
#     field names, and add the new "observaciones", "cita" and attribution rows ---
$ws2.Cells.Item(2, 1).Value = "nomindicador"
$ws2.Cells.Item(2, 2).Value = "Porcentaje de personas que residen en hogares con ingreso bajo la línea de pobreza luego de pagar los gastos de vivienda (alquiler o cuota de compra)"

$ws2.Cells.Item(3, 1).Value = "derecho"
$ws2.Cells.Item(3, 2).Value = "Vivienda"

$ws2.Cells.Item(4, 1).Value = "conindicador"
$ws2.Cells.Item(4, 2).Value = "Gasto excesivo en vivienda en relación a la línea de pobreza"

$ws2.Cells.Item(5, 1).Value = "tipoind"
$ws2.Cells.Item(5, 2).Value = "Resultados"

$ws2.Cells.Item(6, 1).Value = "definicion"
$ws2.Cells.Item(6, 2).Value = "El indicador mide el porcentaje de personas que reside en viviendas con ingreso insuficiente luego del gasto en vivienda. El indicador se construye restando el gasto en vivienda del ingreso de los hogares y comparándolo con el valor de la línea de pobreza a la cual se le quita el peso relativo del componente asociado a la vivienda."

$ws2.Cells.Item(7, 1).Value = "calculo"
$ws2.Cells.Item(7, 2).Value = "Para cada año calcular: (Número de personas en hogares cuyos ingresos totales al deducir los gastos en vivienda por cuota de compra o alquiler caen por debajo de la línea de pobreza quitando el componente de vivienda/ Cantidad total de personas en viviendas particulares)*100"

$ws2.Cells.Item(8, 1).Value = "observaciones"
$ws2.Cells.Item(8, 2).Value = "Sin observaciones"

$ws2.Cells.Item(9, 1).Value = "cita"
$ws2.Cells.Item(9, 2).Value = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE"

$ws2.Cells.Item(10, 1).Value = "Mirador DESCA - UMAD/FCS – INDDHH"
$ws2.Cells.Item(10, 2).Value = " "
